$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing rows 25-27 (content no longer present in the updated sheet)
$ws.Range("A25:A27").EntireRow.Delete()

$ws.Range("B10").Value = "5840747 - Alain Laurent Marie Robin"
$ws.Range("C10").Value = "5840747 - Alain Laurent Marie Robin"

$ws.Range("A13").Value = "Programa resumido:"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Fundamentals of corrosion (thermodynamics and kinetic aspects); Main types of corrosion; Corrosion protection and control; Degradation of polymeric and ceramic materials; Oxidation at high temperatures"
$ws.Range("C14").Value = "Fundamentals of corrosion (thermodynamics and kinetic aspects); Main types of corrosion; Corrosion protection and control; Degradation of polymeric and ceramic materials; Oxidation at high temperatures"

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C15").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Economic, social and environmental effects of metallic corrosion. 2. Thermodynamic and kinetics aspects of corrosion. Polarization. 3.Pourbaix diagram. 4. Types of corrosion (uniform corrosion, pitting corrosion, intergranular corrosion, corrosion associated with mechanical factors, galvanic corrosion, atmospheric corrosion). 5. Anodic and cathodic protection. 6. Coatings. 7. Corrosion inhibitors. 8. Polymeric and ceramic degradation. 9. High temperatures oxidation."
$ws.Range("C16").Value = "1. Economic, social and environmental effects of metallic corrosion. 2. Thermodynamic and kinetics aspects of corrosion. Polarization. 3.Pourbaix diagram. 4. Types of corrosion (uniform corrosion, pitting corrosion, intergranular corrosion, corrosion associated with mechanical factors, galvanic corrosion, atmospheric corrosion). 5. Anodic and cathodic protection. 6. Coatings. 7. Corrosion inhibitors. 8. Polymeric and ceramic degradation. 9. High temperatures oxidation."

$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "As avaliações serão por meio de provas individuais ou trabalhos em equipes, conforme adequação ao conteúdo programático."
$ws.Range("C19").Value = "As avaliações serão por meio de provas individuais ou trabalhos em equipes, conforme adequação ao conteúdo programático."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final é a média aritmética das avaliações realizadas"
$ws.Range("C20").Value = "A nota final é a média aritmética das avaliações realizadas"

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Prova escrita sobre toda matéria.A média final MF será a média da nota final NF e da nota obtida na recuperação NR:MF = (NF + NR)/2Será aprovado o aluno com MF igual ou superior a 5."
$ws.Range("C21").Value = "Prova escrita sobre toda matéria.A média final MF será a média da nota final NF e da nota obtida na recuperação NR:MF = (NF + NR)/2Será aprovado o aluno com MF igual ou superior a 5."

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3003 -  Cinética de Transformação em Materiais  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3003 -  Cinética de Transformação em Materiais  (Requisito fraco)`n"

$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "LOM3008 -  Eletrometalurgia  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOM3008 -  Eletrometalurgia  (Requisito fraco)`n"

# Row height adjustments
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(22).EntireRow.AutoFit()

$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

